$d = $word.ActiveDocument

# 1. Update the ID placeholder text in the first paragraph's first run.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5333_topic_2__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_AFICC_PGI_5333__ID**", 2)

# 2. Remove the trailing space run that followed the ID placeholder.
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute(" ", $true, $false, $false, $false, $false,
                        $true, 1, $false, "", 2)

# 3. Add a paragraph border (5pt space on each side) to the first paragraph.
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5

# 4. Increase the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25
